# Apply Korean translations to the "condition" column (column D) of the kariwaza sheet,
# per the commit "Completed translation of Kariwaza information. Edit typo."
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "집회소★3모래위의 용각(砂上の竜脚)"
$ws.Range("D4").Value = "마을★3귀면수인을 위협하다(鬼面狩人を威す) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D8").Value = "집회소★6유적평원의 흑식룡조사(遺跡平原の黒蝕竜調査)"
$ws.Range("D10").Value = "집회소★5위협! 화산의 철퇴!(脅威！火山の鉄槌！)"
$ws.Range("D13").Value = "마을★4환혹의 마술사(幻惑の魔術師) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D14").Value = "집회소★6사막의 저편에서(砂漠の彼方から)"
$ws.Range("D15").Value = "마을★2도약의 무법자(跳躍のアウトロー) 혹은 집회소★2설산의 주인, 도도블랑고(雪山の主、ドドブランゴ)"
$ws.Range("D17").Value = "집회소★7비리비리바리바리패닉!!(ビリビリバリバリパニック！！)"
$ws.Range("D18").Value = "마을★3귀면수인을 위협하다(鬼面狩人を威す) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D19").Value = "마을★6작열의 칼(灼熱の刃)"
$ws.Range("D20").Value = "집회소★7그대의 힘을, 보여주어보아라(汝ノチカラヲ、見セテミヨ)"
$ws.Range("D21").Value = "마을★7조사대 첫 출전! 유군령의 도모수(調査隊初陣！遺群嶺の桃毛獣) 혹은 집회소★7그 솜씨, 소문대로일까나...?(その腕前、噂通りかしら…？)"
$ws.Range("D22").Value = "G★2겸해의 수렵을 시켜주지!(鎌蟹の狩猟をさせてやるぜ！)"
$ws.Range("D23").Value = "G★4거수의 진격(巨獣の進攻)"
$ws.Range("D25").Value = "마을★4환혹의 마술사(幻惑の魔術師) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D26").Value = "집회소★7고대의 하룡, 오나즈치(古の霞龍、オオナズチ)"
$ws.Range("D27").Value = "마을★3귀면수인을 위협하다(鬼面狩人を威す) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D28").Value = "집회소★3굉룡 티가렉스(轟竜ティガレックス)"
$ws.Range("D29").Value = "집회소★6사냥당하기 전에 사냥해라!(狩られる前に狩れ！)"
$ws.Range("D30").Value = "마을★2도약의 무법자(跳躍のアウトロー) 혹은 집회소★2설산의 주인, 도도블랑고(雪山の主、ドドブランゴ)"
$ws.Range("D32").Value = "집회소★7그랜드・헌터・게임(グランド・ハンター・ゲーム)"
$ws.Range("D33").Value = "마을★7조사대 첫 출전! 유군령의 도모수(調査隊初陣！遺群嶺の桃毛獣) 혹은 집회소★7그 솜씨, 소문대로일까나...?(その腕前、噂通りかしら…？)"
$ws.Range("D34").Value = "G★2하늘의 가족・공로를 지킵니다(空の家族・空路を守るッス)"
$ws.Range("D35").Value = "G★4훈풍의 밀림"
$ws.Range("D36").Value = "마을★3귀면수인을 위협하다(鬼面狩人を威す) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D37").Value = "마을★5괴조의 머리파괴에 도전!(怪鳥の頭部破壊に挑戦！) 혹은 집회소★3쌍두의 뼈(双頭の骸)"
$ws.Range("D38").Value = "집회소★7늪지취몽담(沼地酔夢譚)"
$ws.Range("D40").Value = "마을★4환혹의 마술사(幻惑の魔術師) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D41").Value = "집회소★5땀과 눈물의 연속수렵(汗と涙の連続狩猟)"
$ws.Range("D42").Value = "마을★2도약의 무법자(跳躍のアウトロー) 혹은 집회소★2설산의 주인, 도도블랑고(雪山の主、ドドブランゴ)"
$ws.Range("D44").Value = "집회소★6스릴과 쇼크의 빙해 투어(スリルとショックの氷海ツアー)"
$ws.Range("D45").Value = "마을★7조사대 첫 출전! 유군령의 도모수(調査隊初陣！遺群嶺の桃毛獣) 혹은 집회소★7그 솜씨, 소문대로일까나...?(その腕前、噂通りかしら…？)"
$ws.Range("D46").Value = "G★2밀림의 나르가쿠르가와 대치해라(密林のナルガクルガと対峙せよ)"
$ws.Range("D47").Value = "G★4금과 은이 가져오는 비애(金と銀がもたらす悲哀)"
$ws.Range("D49").Value = "마을★4환혹의 마술사(幻惑の魔術師) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D50").Value = "집회소★5땀과 눈물의 연속수렵(汗と涙の連続狩猟)"
$ws.Range("D51").Value = "마을★2도약의 무법자(跳躍のアウトロー) 혹은 집회소★2설산의 주인, 도도블랑고(雪山の主、ドドブランゴ)"
$ws.Range("D53").Value = "집회소★6광란의 입체투기장(狂乱の立体闘技場)"
$ws.Range("D54").Value = "마을★3귀면수인을 위협하다(鬼面狩人を威す) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D56").Value = "집회소★7분노의 포효(憤怒の雄叫び)"
$ws.Range("D57").Value = "마을★7조사대 첫 출전! 유군령의 도모수(調査隊初陣！遺群嶺の桃毛獣) 혹은 집회소★7그 솜씨, 소문대로일까나...?(その腕前、噂通りかしら…？)"
$ws.Range("D58").Value = "G★2밀림의 나르가쿠르가와 대치해라(密林のナルガクルガと対峙せよ)"
$ws.Range("D59").Value = "G★4초전(焦電)"
$ws.Range("D61").Value = "마을★4환혹의 마술사(幻惑の魔術師) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D62").Value = "집회소★6도볼을 위하여 종은 울리나(ドボルがために銅鑼は鳴る)"
$ws.Range("D63").Value = "마을★2도약의 무법자(跳躍のアウトロー) 혹은 집회소★2설산의 주인, 도도블랑고(雪山の主、ドドブランゴ)"
$ws.Range("D64").Value = "집회소★5빠른 신룡의 수렵피로(疾き迅竜の狩猟披露)"
$ws.Range("D65").Value = "집회소★5그레이트・헌터・게임(グレート・ハンター・ゲーム)"
$ws.Range("D66").Value = "마을★3귀면수인을 위협하다(鬼面狩人を威す) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D67").Value = "마을★5우뚝 솟은 거수(峨々たる巨獣)"
$ws.Range("D68").Value = "집회소★7용맹과감한 브레이브 태스크(勇猛果敢なブレイブタスク)"
$ws.Range("D69").Value = "마을★7조사대 첫 출전! 유군령의 도모수(調査隊初陣！遺群嶺の桃毛獣) 혹은 집회소★7그 솜씨, 소문대로일까나...?(その腕前、噂通りかしら…？)"
$ws.Range("D70").Value = "G★2......원한다면 자격을 보여라(……欲するならば資格を示して)"
$ws.Range("D71").Value = "G★4사막의 사투와 미래(砂漠の死闘と未来)"
$ws.Range("D73").Value = "마을★4환혹의 마술사(幻惑の魔術師) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D74").Value = "집회소★6사냥당하기 전에 사냥해라!(狩られる前に狩れ！)"
$ws.Range("D75").Value = "마을★3귀면수인을 위협하다(鬼面狩人を威す) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D76").Value = "闘技大会★가라라아자라 토벌(ガララアジャラ討伐)"
$ws.Range("D77").Value = "집회소★7시련의 귀결점(試練の帰結点)"
$ws.Range("D78").Value = "마을★2도약의 무법자(跳躍のアウトロー) 혹은 집회소★2설산의 주인, 도도블랑고(雪山の主、ドドブランゴ)"
$ws.Range("D79").Value = "집회소★4대지를 헤엄치는 몬스터(大地を泳ぐモンスター)"
$ws.Range("D80").Value = "집회소★5독, 마비, 혼란에 주의!(毒、麻痺、混乱にご用心！)"
$ws.Range("D81").Value = "마을★7조사대 첫 출전! 유군령의 도모수(調査隊初陣！遺群嶺の桃毛獣) 혹은 집회소★7그 솜씨, 소문대로일까나...?(その腕前、噂通りかしら…？)"
$ws.Range("D82").Value = "G★2......원한다면 자격을 보여라(……欲するならば資格を示して)"
$ws.Range("D83").Value = "G★4사냥꾼의 혼이여 모래속에 잠들어라(狩魂よ砂中に眠れ)"
$ws.Range("D85").Value = "마을★4환혹의 마술사(幻惑の魔術師) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D86").Value = "집회소★5교사룡은 춤추고 노래한다(絞蛇竜は踊り奏でる)"
$ws.Range("D87").Value = "마을★3귀면수인을 위협하다(鬼面狩人を威す) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D89").Value = "집회소★7빙점하의 지배자(氷点下の支配者)"
$ws.Range("D90").Value = "마을★2도약의 무법자(跳躍のアウトロー) 혹은 집회소★2설산의 주인, 도도블랑고(雪山の主、ドドブランゴ)"
$ws.Range("D91").Value = "집회소★3낫장군의 포위진(鎌将軍の包囲陣)"
$ws.Range("D92").Value = "집회소★6원시림의 보통내기가 아닌 자들(原生林の曲者たち)"
$ws.Range("D93").Value = "마을★7조사대 첫 출전! 유군령의 도모수(調査隊初陣！遺群嶺の桃毛獣) 혹은 집회소★7그 솜씨, 소문대로일까나...?(その腕前、噂通りかしら…？)"
$ws.Range("D94").Value = "G★2의뢰위탁 : 화산의 암룡(依頼委託：火山の岩竜)"
$ws.Range("D95").Value = "G★4쉘・위・댄스?(シャル・ウィ・ダンス？)"
$ws.Range("D96").Value = "마을★2도약의 무법자(跳躍のアウトロー) 혹은 집회소★2설산의 주인, 도도블랑고(雪山の主、ドドブランゴ)"
$ws.Range("D97").Value = "집회소★3연옥의 왕, 분개한 염제(煉獄の主、怒れる炎帝)"
$ws.Range("D98").Value = "집회소★7개수일촉의 파워풀 암즈(鎧袖一触のパワフルアームズ)"
$ws.Range("D100").Value = "마을★4환혹의 마술사(幻惑の魔術師) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D101").Value = "집회소★5교사룡은 춤추고 노래한다(絞蛇竜は踊り奏でる)"
$ws.Range("D102").Value = "마을★3귀면수인을 위협하다(鬼面狩人を威す) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D103").Value = "마을★5빨간 해님 아뜨뜨(赤いおひさまアッチッチ)"
$ws.Range("D104").Value = "집회소★7열기로 열광! 불꽃의 군세!(熱気で熱狂！炎の軍勢！)"
$ws.Range("D105").Value = "마을★7조사대 첫 출전! 유군령의 도모수(調査隊初陣！遺群嶺の桃毛獣) 혹은 집회소★7그 솜씨, 소문대로일까나...?(その腕前、噂通りかしら…？)"
$ws.Range("D106").Value = "G★2의뢰위탁 : 화산의 암룡(依頼委託：火山の岩竜)"
$ws.Range("D107").Value = "G★4사신은 갑옷을 두르다(死神は鎧をまとう)"
$ws.Range("D109").Value = "마을★4환혹의 마술사(幻惑の魔術師) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D110").Value = "집회소★6도볼을 위하여 종은 울리나(ドボルがために銅鑼は鳴る)"
$ws.Range("D111").Value = "마을★2도약의 무법자(跳躍のアウトロー) 혹은 집회소★2설산의 주인, 도도블랑고(雪山の主、ドドブランゴ)"
$ws.Range("D112").Value = "집회소★5벽에는 귀가 있고, 천장에는 눈이 있다?(壁に耳あり、天井に目あり？)"
$ws.Range("D113").Value = "집회소★7파괴와 멸망의 산물(破壊と滅亡の申し子)"
$ws.Range("D114").Value = "마을★3귀면수인을 위협하다(鬼面狩人を威す) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D115").Value = "마을★5혼란의 호로로로호루루루(混乱のホロロロホルルル)"
$ws.Range("D116").Value = "집회소★6타버려라! 화산의 열투!!(燃えたぎれ！火山の熱闘！！)"
$ws.Range("D117").Value = "마을★7조사대 첫 출전! 유군령의 도모수(調査隊初陣！遺群嶺の桃毛獣) 혹은 집회소★7그 솜씨, 소문대로일까나...?(その腕前、噂通りかしら…？)"
$ws.Range("D118").Value = "G★2둘이서 특훈, 성과는 반분(二人で特訓、成果は山分け)"
$ws.Range("D119").Value = "G★4기사와 빙해의 결전(騎士と氷海の決闘)"
$ws.Range("D120").Value = "마을★3귀면수인을 위협하다(鬼面狩人を威す) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D122").Value = "집회소★7경천동지할 그랜드 윙(震天動地なグランドウイング)"
$ws.Range("D124").Value = "마을★4환혹의 마술사(幻惑の魔術師) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D125").Value = "집회소★5독, 마비, 혼란에 주의!(毒、麻痺、混乱にご用心！)"
$ws.Range("D126").Value = "마을★2도약의 무법자(跳躍のアウトロー) 혹은 집회소★2설산의 주인, 도도블랑고(雪山の主、ドドブランゴ)"
$ws.Range("D128").Value = "집회소★6비룡들의 난무(飛竜たちの乱舞)"
$ws.Range("D129").Value = "마을★7조사대 첫 출전! 유군령의 도모수(調査隊初陣！遺群嶺の桃毛獣) 혹은 집회소★7그 솜씨, 소문대로일까나...?(その腕前、噂通りかしら…？)"
$ws.Range("D130").Value = "G★2둘이서 특훈, 성과는 반분(二人で特訓、成果は山分け)"
$ws.Range("D131").Value = "G★4무심으로 삼라만상을 자른다(無心にて森羅万象を断つ)"
$ws.Range("D133").Value = "마을★4환혹의 마술사(幻惑の魔術師) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D134").Value = "집회소★6사막의 저편에서(砂漠の彼方から)"
$ws.Range("D135").Value = "마을★2도약의 무법자(跳躍のアウトロー) 혹은 집회소★2설산의 주인, 도도블랑고(雪山の主、ドドブランゴ)"
$ws.Range("D136").Value = "집회소★3사냥당하기 전에 사냥해라!(狩られる前に狩れ！)"
$ws.Range("D137").Value = "집회소★7아직 보지 못한 비탕을 찾아서(まだ見ぬ秘湯をもとめて)"
$ws.Range("D138").Value = "마을★3귀면수인을 위협하다(鬼面狩人を威す) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D139").Value = "집회소★4모이는 강호(集いし強豪)"
$ws.Range("D140").Value = "집회소★6늪지의 광조악단(沼地の狂騒楽団)"
$ws.Range("D141").Value = "마을★7조사대 첫 출전! 유군령의 도모수(調査隊初陣！遺群嶺の桃毛獣) 혹은 집회소★7그 솜씨, 소문대로일까나...?(その腕前、噂通りかしら…？)"
$ws.Range("D142").Value = "G★2겸해의 수렵을 시켜주지!(鎌蟹の狩猟をさせてやるぜ！)"
$ws.Range("D143").Value = "G★4하늘을 돌아, 명계에 내리다(天を廻りて、冥界に堕つ)"
$ws.Range("D144").Value = "마을★3귀면수인을 위협하다(鬼面狩人を威す) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D145").Value = "마을★5취옥의 섬전(翠玉の閃電)"
$ws.Range("D146").Value = "집회소★7숲과언덕의 검은 안개(森丘の黒い霧)"
$ws.Range("D148").Value = "마을★4환혹의 마술사(幻惑の魔術師) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D149").Value = "집회소★5독, 마비, 혼란에 주의!(毒、麻痺、混乱にご用心！)"
$ws.Range("D150").Value = "마을★2도약의 무법자(跳躍のアウトロー) 혹은 집회소★2설산의 주인, 도도블랑고(雪山の主、ドドブランゴ)"
$ws.Range("D151").Value = "집회소★3돌고 울어서 회귀하리라(廻り集いて回帰せん)"
$ws.Range("D152").Value = "집회소★7정들면 고향인 투기장 생활(住めば都の闘技場暮らし)"
$ws.Range("D153").Value = "마을★7조사대 첫 출전! 유군령의 도모수(調査隊初陣！遺群嶺の桃毛獣) 혹은 집회소★7그 솜씨, 소문대로일까나...?(その腕前、噂通りかしら…？)"
$ws.Range("D154").Value = "G★2번뜩임의 한 조각을 찾아서(閃きへの1ピースを求めて)"
$ws.Range("D155").Value = "G★4용암도에서 터지는 쇄광(溶岩島で爆ぜる砕光)"
$ws.Range("D156").Value = "마을★2도약의 무법자(跳躍のアウトロー) 혹은 집회소★2설산의 주인, 도도블랑고(雪山の主、ドドブランゴ)"
$ws.Range("D157").Value = "마을★6분골쇄룡!(粉骨砕竜！)"
$ws.Range("D158").Value = "집회소★7연옥의 왕, 분개한 염제(煉獄の主、怒れる炎帝)"
$ws.Range("D159").Value = "마을★3귀면수인을 위협하다(鬼面狩人を威す) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D160").Value = "마을★5담홍의 포호가 흔들거리는가(淡紅の泡狐がたゆたうか)"
$ws.Range("D161").Value = "집회소★5불면증인 당신에게 최면요법×2(不眠のあなたに催眠療法×２)"
$ws.Range("D163").Value = "마을★4환혹의 마술사(幻惑の魔術師) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D164").Value = "집회소★6지저화산에 울려퍼지는 침략의 발소리(地底火山に響く侵略の足音)"
$ws.Range("D165").Value = "마을★7조사대 첫 출전! 유군령의 도모수(調査隊初陣！遺群嶺の桃毛獣) 혹은 집회소★7그 솜씨, 소문대로일까나...?(その腕前、噂通りかしら…？)"
$ws.Range("D166").Value = "G★2하늘의 가족・공로를 지킵니다(空の家族・空路を守るッス)"
$ws.Range("D167").Value = "G★4전율의 유군령(戦慄の遺群嶺)"
$ws.Range("D169").Value = "마을★4환혹의 마술사(幻惑の魔術師) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D170").Value = "집회소★7기기괴괴의 하드비크(奇奇怪怪のハードビーク)"
$ws.Range("D171").Value = "마을★3귀면수인을 위협하다(鬼面狩人を威す) 혹은 집회소★3용의 반역자(電の反逆者)"
$ws.Range("D172").Value = "투기대회★케챠와챠 토벌(ケチャワチャ討伐)"
$ws.Range("D173").Value = "집회소★5불면증인 당신에게 최면요법×2(不眠のあなたに催眠療法×２)"
$ws.Range("D174").Value = "마을★2도약의 무법자(跳躍のアウトロー) 혹은 집회소★2설산의 주인, 도도블랑고(雪山の主、ドドブランゴ)"
$ws.Range("D175").Value = "마을★6파랑과 녹색의 파상포위망(青と緑の波状包囲網) 혹은 집회소★5계류의 수룡(渓流の水竜)"
$ws.Range("D176").Value = "집회소★6와일드 배럿(ワイルドバレット)"
$ws.Range("D177").Value = "마을★7조사대 첫 출전! 유군령의 도모수(調査隊初陣！遺群嶺の桃毛獣) 혹은 집회소★7그 솜씨, 소문대로일까나...?(その腕前、噂通りかしら…？)"
$ws.Range("D178").Value = "G★2번뜩임의 한 조각을 찾아서(閃きへの1ピースを求めて)"
$ws.Range("D179").Value = "G★4금빛초거성(金色超巨星)"

# Column D is now much wider after translation; drop the old bestFit and widen it
# (target stored width 129.7109375 characters; this host quantizes ColumnWidth to
# 1/6-character steps, so 128.8 is the closest input that lands on the nearest step).
$ws.Columns.Item(4).ColumnWidth = 128.8

# Restore the author's last scroll position/selection in the sheet view
$win = $excel.ActiveWindow
$win.ScrollRow = 166
$win.ScrollColumn = 1
$ws.Range("D180").Select()

